$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New data rows (9, 10, 11) - three more LED/indicator circuits
#    (IMU, Voltage Sensor, Driver circuits per commit message)
# ---------------------------------------------------------------------------

# -- Copy cell formatting (styles) from the existing template row (row 4) so
#    the new rows reuse the same cellXfs entries instead of creating new ones.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C9:C11").PasteSpecial(-4122) | Out-Null

$ws.Range("D4").Copy() | Out-Null
$ws.Range("D9:D11").PasteSpecial(-4122) | Out-Null

$ws.Range("E4").Copy() | Out-Null
$ws.Range("E9:E11").PasteSpecial(-4122) | Out-Null

$ws.Range("F4").Copy() | Out-Null
$ws.Range("F9:F11").PasteSpecial(-4122) | Out-Null

$ws.Range("G4").Copy() | Out-Null
$ws.Range("G9:G11").PasteSpecial(-4122) | Out-Null

$ws.Range("H4").Copy() | Out-Null
$ws.Range("H9:H11").PasteSpecial(-4122) | Out-Null

$ws.Range("I4").Copy() | Out-Null
$ws.Range("I9:I11").PasteSpecial(-4122) | Out-Null

# J9:J11 use the centred style (same as C4/D4 ...), not the J4 style.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("J9:J11").PasteSpecial(-4122) | Out-Null

$ws.Range("K4").Copy() | Out-Null
$ws.Range("K9:K11").PasteSpecial(-4122) | Out-Null

$ws.Range("L4").Copy() | Out-Null
$ws.Range("L9:L11").PasteSpecial(-4122) | Out-Null

# -- Values (written in the exact order the strings were first introduced so
#    the shared-string table ends up in the same order).
$ws.Range("A9").Value = "JLCPCB"

$ws.Range("E11").Value = "Yellow"
$ws.Range("E9").Value = "Green"

$ws.Range("B9").Value = "C2297"
$ws.Range("B10").Value = "C2286"
$ws.Range("B11").Value = "C2296"

$ws.Range("D9").Value = "SMD"
$ws.Range("D10").Value = "SMD"
$ws.Range("D11").Value = "SMD"

$ws.Range("E10").Value = "RED"

$ws.Range("C9").Value = 5
$ws.Range("F9").Value = 2.8
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 5
$ws.Range("J9").Value = 470

$ws.Range("C10").Value = 3.3
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = 20
$ws.Range("J10").Value = 75

$ws.Range("C11").Value = 3.3
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = 20
$ws.Range("J11").Value = 75

# -- Formulas: I column stays per-cell (not shared), K/L are shared across
#    the new block, matching how the original table was built.
$ws.Range("I9").Formula = "=(C9-F9)/(H9*10^-3)"
$ws.Range("I10").Formula = "=(C10-F10)/(H10*10^-3)"
$ws.Range("I11").Formula = "=(C11-F11)/(H11*10^-3)"

$ws.Range("K9:K11").Formula = "=(C9-F9)/J9*10^3"
$ws.Range("L9:L11").Formula = "=(C9-F9)*K9*10^-3"

# ---------------------------------------------------------------------------
# 2. Conditional formatting: keep the existing "red" rule on L4:L7, and add a
#    duplicate rule (same look) on the new L9:L11 range with top priority.
# ---------------------------------------------------------------------------
$oldRule = $ws.Range("L4:L7").FormatConditions.Item(1)
$oldRule.Priority = 4

$newRule = $ws.Range("L9:L11").FormatConditions.Add(1, 5, "0.25")
$newRule.Font.Color = 393372
$newRule.Interior.Color = 13551615
$newRule.Priority = 1

# ---------------------------------------------------------------------------
# 3. Reposition the two pictures that moved further down/right on the sheet
#    to make room for the new rows.
# ---------------------------------------------------------------------------
$bigPic = $ws.Shapes.Item(2)
$bigPic.Left = 75.58692913385828
$bigPic.Top = 184.10874015748033

$smallPic = $ws.Shapes.Item(3)
$smallPic.Left = 840.7500787401575
$smallPic.Top = 132.75
$smallPic.Width = 47.63826771653544
$smallPic.Height = 184.55141732283465

# ---------------------------------------------------------------------------
# 4. Selection moves to J9.
# ---------------------------------------------------------------------------
$ws.Range("J9").Select() | Out-Null
